$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell G2 value
$ws.Range("G2").Value = '${NULL}'

# Add new column H with header (row1) and value (row2)
$ws.Range("H1").Value = "Caso6"
$ws.Range("H2").Value = '${EMPTY}'

# Update the selection to match the target state
$ws.Range("I7:J8").Select()
